$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update changed cell values
$ws.Range("F9").Value = 22.35
$ws.Range("D21").Value = 11.04833333333333
$ws.Range("F21").Value = 15.9
$ws.Range("F22").Value = 16.5
$ws.Range("D25").Value = 8.318333333333333
$ws.Range("F25").Value = 14.9
$ws.Range("D29").Value = 11.78333333333333
$ws.Range("F29").Value = 16.65
$ws.Range("D33").Value = 8.773333333333333
$ws.Range("F33").Value = 15.15
$ws.Range("D41").Value = 9.263333333333332
$ws.Range("F41").Value = 13.15
$ws.Range("F42").Value = 14.9
$ws.Range("D43").Value = 10.31333333333333
$ws.Range("F43").Value = 13.95
$ws.Range("D49").Value = 10.06833333333333
$ws.Range("F49").Value = 15.4
$ws.Range("D50").Value = 8.668333333333333
$ws.Range("F50").Value = 13.6
$ws.Range("D52").Value = 12.90333333333333
$ws.Range("F52").Value = 14.65
$ws.Range("D59").Value = 9.893333333333333
$ws.Range("F59").Value = 11.95
$ws.Range("D60").Value = 13.88333333333333
$ws.Range("F60").Value = 19.55
$ws.Range("D68").Value = 12.03416666666667
$ws.Range("F68").Value = 16.1
$ws.Range("F73").Value = 13.8
$ws.Range("D81").Value = 10.63416666666667
$ws.Range("F81").Value = 14.7
$ws.Range("D86").Value = 9.648333333333332
$ws.Range("F86").Value = 18.2
$ws.Range("F93").Value = 13.2
$ws.Range("D94").Value = 8.423333333333332
$ws.Range("F94").Value = 14.25
$ws.Range("D100").Value = 10.52916666666667
$ws.Range("F100").Value = 13.6
$ws.Range("F103").Value = 18.1
$ws.Range("D107").Value = 10.13483333333333
$ws.Range("F107").Value = 14.6
$ws.Range("D108").Value = 6.704833333333333
$ws.Range("F108").Value = 10.05
$ws.Range("D119").Value = 8.744166666666665
$ws.Range("F119").Value = 13.1
$ws.Range("D124").Value = 7.478333333333333
$ws.Range("F124").Value = 8.5
$ws.Range("F125").Value = 8.5
$ws.Range("F127").Value = 9.050000000000001
$ws.Range("F130").Value = 8.366666666666667
$ws.Range("F137").Value = 10.8
$ws.Range("F138").Value = 7.8
$ws.Range("D141").Value = 2.473333333333332
$ws.Range("F141").Value = 7.75
$ws.Range("F145").Value = 11.03333333333333
$ws.Range("D147").Value = 4.118333333333334
$ws.Range("F147").Value = 4.9
$ws.Range("F149").Value = 11.4
$ws.Range("F150").Value = 1.7
$ws.Range("F162").Value = 2.3
$ws.Range("D172").Value = 5.168333333333333
$ws.Range("F172").Value = 7
$ws.Range("D193").Value = 2.718333333333331
$ws.Range("F193").Value = 8.6
$ws.Range("F195").Value = 9.725
$ws.Range("D210").Value = 5.066833333333332
$ws.Range("F210").Value = 6.949999999999999
$ws.Range("F216").Value = 6.875
$ws.Range("D226").Value = 5.178833333333333
$ws.Range("F226").Value = 7.3
$ws.Range("F230").Value = 10
$ws.Range("F239").Value = 6.3
$ws.Range("F241").Value = 0.8
$ws.Range("F247").Value = 3.85
$ws.Range("D248").Value = 1.937833333333333
$ws.Range("F248").Value = 2.3
$ws.Range("D250").Value = 2.823333333333332
$ws.Range("F250").Value = 5.25
$ws.Range("D251").Value = 4.188333333333332
$ws.Range("F251").Value = 12.2
$ws.Range("F252").Value = 13.45
$ws.Range("F253").Value = 5.1
$ws.Range("D256").Value = 11.67833333333333
$ws.Range("F256").Value = 13.8
$ws.Range("F257").Value = 5.699999999999999
$ws.Range("F261").Value = 1.6
$ws.Range("D265").Value = 5.553333333333332
$ws.Range("F265").Value = 4.65
$ws.Range("F274").Value = 6.066666666666667
$ws.Range("F276").Value = 6.9
$ws.Range("F285").Value = 14.9
$ws.Range("F287").Value = 5
$ws.Range("F290").Value = 5
$ws.Range("D293").Value = 4.573333333333333
$ws.Range("F293").Value = 6.35
$ws.Range("D307").Value = 6.171666666666666
$ws.Range("F307").Value = 15.03333333333333
$ws.Range("F308").Value = 1.8
$ws.Range("F309").Value = 3
$ws.Range("F311").Value = 2.3
$ws.Range("D315").Value = 0.548333333333332
$ws.Range("F315").Value = 1
$ws.Range("F317").Value = 2.3
$ws.Range("F320").Value = 3
$ws.Range("F324").Value = 10.36666666666667
$ws.Range("F327").Value = 3.5
$ws.Range("D328").Value = 1.458333333333334
$ws.Range("F328").Value = 5.1
$ws.Range("D341").Value = 8.306666666666667
$ws.Range("F341").Value = 11.475
$ws.Range("D361").Value = -1.551666666666668
$ws.Range("F361").Value = 2
$ws.Range("D365").Value = 0.0688333333333334
$ws.Range("F375").Value = 1.5
$ws.Range("F380").Value = 1.4
$ws.Range("F385").Value = 5.433333333333334
$ws.Range("D388").Value = 0.0583333333333325
$ws.Range("F389").Value = 3.766666666666667
$ws.Range("F427").Value = 2.133333333333333
$ws.Range("D428").Value = 4.299166666666665
$ws.Range("F428").Value = 6.75
$ws.Range("D434").Value = 0.4421666666666668
$ws.Range("F434").Value = 0.5333333333333333
$ws.Range("F439").Value = 1
$ws.Range("D456").Value = 2.298333333333332
$ws.Range("F490").Value = 2
$ws.Range("F491").Value = 0.2
$ws.Range("D493").Value = -0.3301666666666668
$ws.Range("F501").Value = 3
$ws.Range("D506").Value = -3.091666666666667
$ws.Range("D511").Value = 0.548333333333332
$ws.Range("F523").Value = 0.325
$ws.Range("D525").Value = -2.951666666666667
$ws.Range("F525").Value = 0.1
$ws.Range("F532").Value = 1.6
$ws.Range("F536").Value = 0.8
$ws.Range("D543").Value = 0.8983333333333324
$ws.Range("F543").Value = 1.2
$ws.Range("D561").Value = 0.00583333333333325
$ws.Range("D678").Value = -3.651666666666668
$ws.Range("F678").Value = -2

# Clear cells that are removed entirely in the target (empty cell, no value)
$ws.Range("F514").ClearContents()
$ws.Range("F556").ClearContents()
$ws.Range("F618").ClearContents()
$ws.Range("F654").ClearContents()
$ws.Range("F667").ClearContents()

Write-Output "Applied all cell updates"
